$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "SourceCode" styled paragraph right after the Date
#    paragraph ("February 16, 2017"), containing the first R chunk:
#      knitr::opts_chunk$set(echo = TRUE)
#      <blank line>
#      # define libraries
#      library(ggplot2)
# ------------------------------------------------------------------

$datePara = $d.Paragraphs.Item(3)
$afterDate = $datePara.Range
$afterDate.Collapse(0)
$afterDate.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(4)
$newPara.Style = "SourceCode"

$insertRng = $newPara.Range
$insertRng.Collapse(1)

# BR marks a soft line-break (w:br) rather than literal text.
$BR = [char]11

$segments = @(
    @{Text="knitr::opts_chunk$"; Style="NormalTok"},
    @{Text="set"; Style="KeywordTok"},
    @{Text="("; Style="NormalTok"},
    @{Text="echo ="; Style="DataTypeTok"},
    @{Text=" "; Style="NormalTok"},
    @{Text="TRUE"; Style="OtherTok"},
    @{Text=")"; Style="NormalTok"},
    @{Text=$BR; Style=$null},
    @{Text=$BR; Style=$null},
    @{Text="# define libraries"; Style="CommentTok"},
    @{Text=$BR; Style=$null},
    @{Text="library"; Style="KeywordTok"},
    @{Text="(ggplot2)"; Style="NormalTok"}
)

foreach ($seg in $segments) {
    $insertRng.Collapse(0)
    $segStart = $insertRng.Start
    $insertRng.InsertAfter($seg.Text)
    $segEnd = $insertRng.End
    if ($seg.Style) {
        $segRange = $d.Range($segStart, $segEnd)
        $segRange.Style = $seg.Style
    }
}

# NOTE: the source diff also shows two <w:nsid> GUIDs changing inside
# numbering.xml's <w:abstractNum> definitions. That nsid value is an
# internal list-numbering identifier that Word never surfaces through
# the VBA/COM object model (no List/ListTemplate/ListGallery property
# maps to it) -- it is simply re-minted by the authoring pipeline
# (pandoc/knitr) each time the document is regenerated from the Rmd
# source, independent of any in-document edit. There is no COM-visible
# API to target it, so it is intentionally left untouched here.
